$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.150782
$ws.Range("H2").Value = 3.452345999999999
$ws.Range("I2").Value = 0.03823856951930295
$ws.Range("J2").Value = 0.03823856951930295
$ws.Range("M2").Value = 45.90594266666667
$ws.Range("N2").Value = 137.717828
$ws.Range("O2").Value = 0.3954672001633582
$ws.Range("P2").Value = 0.3954672001633583
$ws.Range("Q2").Value = 52.82773251383199
$ws.Range("R2").Value = 475.4495926244879
$ws.Range("S2").Value = 0.01512210002605067
$ws.Range("T2").Value = 0.01512210002605067
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.150782
$ws.Range("H3").Value = 3.452345999999999
$ws.Range("I3").Value = 0.03823856951930295
$ws.Range("J3").Value = 0.03823856951930295
$ws.Range("O3").Value = 0.3484294080560655
$ws.Range("P3").Value = 0.3484294080560656
$ws.Range("Q3").Value = 46.54427866871199
$ws.Range("R3").Value = 418.8985080184079
$ws.Range("S3").Value = 0.01332344214252144
$ws.Range("T3").Value = 0.01332344214252144
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.150782
$ws.Range("H4").Value = 3.452345999999999
$ws.Range("I4").Value = 0.03823856951930295
$ws.Range("J4").Value = 0.03823856951930295
$ws.Range("M4").Value = 12.761795
$ws.Range("N4").Value = 38.28538500000001
$ws.Range("O4").Value = 0.1099393900775594
$ws.Range("P4").Value = 0.1099393900775594
$ws.Range("Q4").Value = 14.68604397369
$ws.Range("R4").Value = 132.17439576321
$ws.Range("S4").Value = 0.004203925010390521
$ws.Range("T4").Value = 0.004203925010390521
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.150782
$ws.Range("H5").Value = 3.452345999999999
$ws.Range("I5").Value = 0.03823856951930295
$ws.Range("J5").Value = 0.03823856951930295
$ws.Range("M5").Value = 16.966758
$ws.Range("N5").Value = 50.900274
$ws.Range("O5").Value = 0.1461640017030168
$ws.Range("P5").Value = 0.1461640017030168
$ws.Range("Q5").Value = 19.52503970475599
$ws.Range("R5").Value = 175.725357342804
$ws.Range("S5").Value = 0.005589102340340323
$ws.Range("T5").Value = 0.005589102340340323
$ws.Range("I6").Value = 0.9169230158851821
$ws.Range("J6").Value = 0.916923015885182
$ws.Range("M6").Value = 45.90594266666667
$ws.Range("N6").Value = 137.717828
$ws.Range("O6").Value = 0.3954672001633582
$ws.Range("P6").Value = 0.3954672001633583
$ws.Range("Q6").Value = 1266.756691682893
$ws.Range("R6").Value = 11400.81022514604
$ws.Range("S6").Value = 0.3626129778574554
$ws.Range("T6").Value = 0.3626129778574554
$ws.Range("I7").Value = 0.9169230158851821
$ws.Range("J7").Value = 0.916923015885182
$ws.Range("O7").Value = 0.3484294080560655
$ws.Range("P7").Value = 0.3484294080560656
$ws.Range("S7").Value = 0.3194829436578564
$ws.Range("T7").Value = 0.3194829436578564
$ws.Range("I8").Value = 0.9169230158851821
$ws.Range("J8").Value = 0.916923015885182
$ws.Range("M8").Value = 12.761795
$ws.Range("N8").Value = 38.28538500000001
$ws.Range("O8").Value = 0.1099393900775594
$ws.Range("P8").Value = 0.1099393900775594
$ws.Range("Q8").Value = 352.1567857024717
$ws.Range("R8").Value = 3169.411071322246
$ws.Range("S8").Value = 0.1008059571144932
$ws.Range("T8").Value = 0.1008059571144932
$ws.Range("I9").Value = 0.9169230158851821
$ws.Range("J9").Value = 0.916923015885182
$ws.Range("M9").Value = 16.966758
$ws.Range("N9").Value = 50.900274
$ws.Range("O9").Value = 0.1461640017030168
$ws.Range("P9").Value = 0.1461640017030168
$ws.Range("Q9").Value = 468.191109563482
$ws.Range("R9").Value = 4213.719986071338
$ws.Range("S9").Value = 0.1340211372553771
$ws.Range("T9").Value = 0.1340211372553771
$ws.Range("G10").Value = 1.290098666666667
$ws.Range("H10").Value = 3.870296
$ws.Range("I10").Value = 0.04286783035543951
$ws.Range("J10").Value = 0.0428678303554395
$ws.Range("M10").Value = 45.90594266666667
$ws.Range("N10").Value = 137.717828
$ws.Range("O10").Value = 0.3954672001633582
$ws.Range("P10").Value = 0.3954672001633583
$ws.Range("Q10").Value = 59.2231954263431
$ws.Range("R10").Value = 533.008758837088
$ws.Range("S10").Value = 0.01695282084774348
$ws.Range("T10").Value = 0.01695282084774348
$ws.Range("G11").Value = 1.290098666666667
$ws.Range("H11").Value = 3.870296
$ws.Range("I11").Value = 0.04286783035543951
$ws.Range("J11").Value = 0.0428678303554395
$ws.Range("O11").Value = 0.3484294080560655
$ws.Range("P11").Value = 0.3484294080560656
$ws.Range("Q11").Value = 52.17905029055644
$ws.Range("R11").Value = 469.6114526150079
$ws.Range("S11").Value = 0.01493641275539363
$ws.Range("T11").Value = 0.01493641275539362
$ws.Range("G12").Value = 1.290098666666667
$ws.Range("H12").Value = 3.870296
$ws.Range("I12").Value = 0.04286783035543951
$ws.Range("J12").Value = 0.0428678303554395
$ws.Range("M12").Value = 12.761795
$ws.Range("N12").Value = 38.28538500000001
$ws.Range("O12").Value = 0.1099393900775594
$ws.Range("P12").Value = 0.1099393900775594
$ws.Range("Q12").Value = 16.46397471377333
$ws.Range("R12").Value = 148.17577242396
$ws.Range("S12").Value = 0.004712863123225306
$ws.Range("T12").Value = 0.004712863123225306
$ws.Range("G13").Value = 1.290098666666667
$ws.Range("H13").Value = 3.870296
$ws.Range("I13").Value = 0.04286783035543951
$ws.Range("J13").Value = 0.0428678303554395
$ws.Range("M13").Value = 16.966758
$ws.Range("N13").Value = 50.900274
$ws.Range("O13").Value = 0.1461640017030168
$ws.Range("P13").Value = 0.1461640017030168
$ws.Range("Q13").Value = 21.888791873456
$ws.Range("R13").Value = 196.999126861104
$ws.Range("S13").Value = 0.006265733629077095
$ws.Range("T13").Value = 0.006265733629077096
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.05930433333333333
$ws.Range("H14").Value = 0.177913
$ws.Range("I14").Value = 0.001970584240075516
$ws.Range("J14").Value = 0.001970584240075516
$ws.Range("M14").Value = 45.90594266666667
$ws.Range("N14").Value = 137.717828
$ws.Range("O14").Value = 0.3954672001633582
$ws.Range("P14").Value = 0.3954672001633583
$ws.Range("Q14").Value = 2.722421325884889
$ws.Range("R14").Value = 24.501791932964
$ws.Range("S14").Value = 0.0007793014321087031
$ws.Range("T14").Value = 0.0007793014321087031
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.05930433333333333
$ws.Range("H15").Value = 0.177913
$ws.Range("I15").Value = 0.001970584240075516
$ws.Range("J15").Value = 0.001970584240075516
$ws.Range("O15").Value = 0.3484294080560655
$ws.Range("P15").Value = 0.3484294080560656
$ws.Range("Q15").Value = 2.398610177191555
$ws.Range("R15").Value = 21.587491594724
$ws.Range("S15").Value = 0.0006866095002941238
$ws.Range("T15").Value = 0.0006866095002941238
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.05930433333333333
$ws.Range("H16").Value = 0.177913
$ws.Range("I16").Value = 0.001970584240075516
$ws.Range("J16").Value = 0.001970584240075516
$ws.Range("M16").Value = 12.761795
$ws.Range("N16").Value = 38.28538500000001
$ws.Range("O16").Value = 0.1099393900775594
$ws.Range("P16").Value = 0.1099393900775594
$ws.Range("Q16").Value = 0.7568297446116666
$ws.Range("R16").Value = 6.811467701505
$ws.Range("S16").Value = 0.0002166448294503531
$ws.Range("T16").Value = 0.0002166448294503531
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.05930433333333333
$ws.Range("H17").Value = 0.177913
$ws.Range("I17").Value = 0.001970584240075516
$ws.Range("J17").Value = 0.001970584240075516
$ws.Range("M17").Value = 16.966758
$ws.Range("N17").Value = 50.900274
$ws.Range("O17").Value = 0.1461640017030168
$ws.Range("P17").Value = 0.1461640017030168
$ws.Range("Q17").Value = 1.006202272018
$ws.Range("R17").Value = 9.055820448161999
$ws.Range("S17").Value = 0.0002880284782223358
$ws.Range("T17").Value = 0.0002880284782223358
